$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header summary values ---
$ws.Range("E11").Value = 341640      # VALOR MORA total
$ws.Range("C13").Value = 1           # Cant. Trabajadores
$ws.Range("F13").Value = 6           # Cant. Periodos

# --- Copy the "last data row" formatting (bottom border) from row 26 down to row 21,
#     which will become the new last data row once the extra rows are removed.
$ws.Range("B26:J26").Copy()
$ws.Range("B21:J21").PasteSpecial(-4122)   # xlPasteFormats

# --- Replace worker data: single worker (ALBERTO JOSE BLANCO JULIO) across periods 2503-2508
$periods = @("2503", "2504", "2505", "2506", "2507", "2508")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Cells.Item($r, 2).Value = "CC"
    $ws.Cells.Item($r, 3).Value = "92400506"
    $ws.Cells.Item($r, 4).Value = "ALBERTO JOSE BLANCO JULIO"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 56940
    $ws.Cells.Item($r, 7).Value = 1423500
}

# --- Remove the now-unused trailing data rows (old rows 22-26) ---
$ws.Rows("22:26").Delete()

Write-Output "done"
